$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'22.328.70"
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = "'1.564.85"
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'1.001"
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = "'286.76"
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').Value = "'0.3772"
$ws.Range('E7').Value = '  +3.09%  '
$ws.Range('D8').Value = "'0.3269"
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('E9').Value = '  -5.42%  '
$ws.Range('D10').Value = "'1.141"
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').Value = "'0.07409"
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').Value = "'1.002"
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = "'20.38"
$ws.Range('E13').Value = '  -2.51%  '
$ws.Range('D14').Value = "'5.844"
$ws.Range('E14').Value = '  -2.58%  '
$ws.Range('D15').Value = "'6.802"
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').Value = "'1.580.82"
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = "'0.06721"
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').Value = "'85.92"
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('D20').Value = "'1.001"
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = "'6.356"
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('D22').Value = "'16.26"
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').Value = "'11.68"
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').Value = "'22.351.32"
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').Value = "'2.300"
$ws.Range('E25').Value = '  -3.70%  '
$ws.Range('E26').Value = '  -4.04%  '
$ws.Range('D27').Value = "'150.51"
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('D28').Value = "'19.39"
$ws.Range('E28').Value = '  -1.59%  '
$ws.Range('D29').Value = "'4.900"
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('D30').Value = "'123.15"
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('D31').Value = "'1.753.07"
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = "'1.048"
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').Value = "'5.920"
$ws.Range('E33').Value = '  -4.25%  '
$ws.Range('D34').Value = "'1.914"
$ws.Range('E34').Value = '  -4.07%  '
$ws.Range('D35').Value = "'9.466"
$ws.Range('E35').Value = '  -3.52%  '
$ws.Range('D36').Value = "'0.08250"
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('D38').Value = "'0.06295"
$ws.Range('E38').Value = '  -3.06%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = "'0.2180"
$ws.Range('E39').Value = '  -3.85%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = "'1.275"
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').Value = "'5.258"
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').Value = "'11.01"
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').Value = "'0.6077"
$ws.Range('E43').Value = '  -3.35%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('E45').Value = '  -1.81%  '
$ws.Range('D46').Value = "'3.746"
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').Value = "'0.5889"
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'1.996"
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = "'123.98"
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').Value = "'1.179"
$ws.Range('E50').Value = '  -3.64%  '
$ws.Range('D51').Value = "'0.07125"
$ws.Range('E51').Value = '  -1.47%  '
